$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the last genre entry from "Mainstream" to "Drum & Bass"
$ws.Range("A13").Value = "Drum & Bass"

# Move the active selection to A13 (single cell), matching the commit's
# "click" stability fix which re-targets the selector.
$ws.Range("A13").Select()
